$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 'Knives Out (2019)'
$ws.Range("B5").Value = '''92%'
$ws.Range("C5").Value = 'Drama, Mystery & Suspense'
$ws.Range("A6").Value = 'Toy Story 4 (2019)'
$ws.Range("B6").Value = '''94%'
$ws.Range("C6").Value = 'Animation, Comedy, Kids & Family, Science Fiction & Fantasy'
$ws.Range("A8").Value = 'Little Women (2019)'
$ws.Range("B8").Value = '''92%'
$ws.Range("C8").Value = 'Drama'
$ws.Range("A10").Value = 'The Farewell (2019)'
$ws.Range("B10").Value = '''87%'
$ws.Range("C10").Value = 'Comedy, Drama'
$ws.Range("A12").Value = 'A Beautiful Day in the Neighborhood (2019)'
$ws.Range("B12").Value = '''92%'
$ws.Range("C12").Value = 'Drama'
$ws.Range("A13").Value = 'Spider-Man: Far From Home (2019)'
$ws.Range("B13").Value = '''95%'
$ws.Range("C13").Value = 'Action & Adventure, Science Fiction & Fantasy'
$ws.Range("A14").Value = 'If Beale Street Could Talk (2019)'
$ws.Range("B14").Value = '''70%'
$ws.Range("C14").Value = 'Drama, Romance'
$ws.Range("B17").Value = '''82%'
$ws.Range("B25").Value = '''81%'
$ws.Range("B26").Value = '''66%'
$ws.Range("B38").Value = '''87%'
$ws.Range("A40").Value = 'Hustlers (2019)'
$ws.Range("B40").Value = '''65%'
$ws.Range("C40").Value = 'Drama'
$ws.Range("A41").Value = 'Atlantics (2019)'
$ws.Range("B41").Value = '''81%'
$ws.Range("C41").Value = 'Drama, Romance'
$ws.Range("A42").Value = 'One Child Nation (2019)'
$ws.Range("B42").Value = '''84%'
$ws.Range("C42").Value = 'Documentary'
$ws.Range("A43").Value = '3 Faces (2019)'
$ws.Range("B43").Value = '''67%'
$ws.Range("C43").Value = 'Art House & International, Drama'
$ws.Range("B45").Value = '''89%'
$ws.Range("A46").Value = 'Chained for Life (2019)'
$ws.Range("B46").Value = '''63%'
$ws.Range("C46").Value = 'Drama'
$ws.Range("A47").Value = 'Captain Marvel (2019)'
$ws.Range("B47").Value = '''48%'
$ws.Range("C47").Value = 'Action & Adventure, Science Fiction & Fantasy'
$ws.Range("A48").Value = 'Midnight Family (2019)'
$ws.Range("B48").Value = '''78%'
$ws.Range("C48").Value = 'Documentary'
$ws.Range("A49").Value = 'First Love (Hatsukoi) (2019)'
$ws.Range("B49").Value = '''88%'
$ws.Range("C49").Value = 'Action & Adventure, Drama, Mystery & Suspense'
$ws.Range("A50").Value = 'The Chambermaid (La camarista) (2019)'
$ws.Range("B50").Value = '''60%'
$ws.Range("C50").Value = 'Art House & International, Drama'
$ws.Range("A51").Value = 'Deadwood: The Movie (2019)'
$ws.Range("B51").Value = '''79%'
$ws.Range("C51").Value = 'Drama, Western'
$ws.Range("A52").Value = 'Midnight Traveler (2019)'
$ws.Range("B52").Value = '''84%'
$ws.Range("C52").Value = 'Documentary'
$ws.Range("A53").Value = 'Transit (2019)'
$ws.Range("B53").Value = '''64%'
$ws.Range("C53").Value = 'Art House & International, Drama'
$ws.Range("A54").Value = 'Mickey and the Bear (2019)'
$ws.Range("B54").Value = '''79%'
$ws.Range("C54").Value = 'Drama'
$ws.Range("A55").Value = 'Tigers Are Not Afraid (Vuelven) (2019)'
$ws.Range("B55").Value = '''86%'
$ws.Range("C55").Value = 'Art House & International, Drama, Horror, Science Fiction & Fantasy'
$ws.Range("A56").Value = 'Hail Satan? (2019)'
$ws.Range("B56").Value = '''85%'
$ws.Range("A57").Value = 'Toni Morrison: The Pieces I Am (2019)'
$ws.Range("B57").Value = '''96%'
$ws.Range("C57").Value = 'Documentary'
$ws.Range("A58").Value = 'Ready or Not (2019)'
$ws.Range("B58").Value = '''78%'
$ws.Range("C58").Value = 'Horror, Mystery & Suspense'
$ws.Range("A59").Value = 'Ad Astra (2019)'
$ws.Range("B59").Value = '''40%'
$ws.Range("C59").Value = 'Mystery & Suspense, Science Fiction & Fantasy'
$ws.Range("A60").Value = 'Midsommar (2019)'
$ws.Range("B60").Value = '''63%'
$ws.Range("C60").Value = 'Horror'
$ws.Range("A61").Value = 'The Heiresses (Las Herederas) (2019)'
$ws.Range("B61").Value = '''70%'
$ws.Range("C61").Value = 'Drama'
$ws.Range("A62").Value = 'Wild Rose (2019)'
$ws.Range("B62").Value = '''87%'
$ws.Range("A63").Value = 'Blinded by the Light (2019)'
$ws.Range("B63").Value = '''91%'
$ws.Range("C63").Value = 'Comedy, Drama'
$ws.Range("A64").Value = '63 Up (2019)'
$ws.Range("B64").Value = '''93%'
$ws.Range("C64").Value = 'Documentary'
$ws.Range("A65").Value = 'Too Late to Die Young (Tarde para morir joven) (2019)'
$ws.Range("B65").Value = '''54%'
$ws.Range("C65").Value = 'Art House & International, Drama'
$ws.Range("A66").Value = 'Homecoming: A Film by Beyoncé (2019)'
$ws.Range("B66").Value = '''87%'
$ws.Range("C66").Value = 'Documentary, Musical & Performing Arts'
$ws.Range("A67").Value = 'By the Grace of God (Grâce à Dieu) (2019)'
$ws.Range("B67").Value = '''77%'
$ws.Range("C67").Value = 'Drama'
$ws.Range("A68").Value = 'Little Woods (2019)'
$ws.Range("B68").Value = '''54%'
$ws.Range("A69").Value = 'Rojo (2019)'
$ws.Range("B69").Value = '''48%'
$ws.Range("C69").Value = 'Drama, Mystery & Suspense'
$ws.Range("A70").Value = 'Gloria Bell (2019)'
$ws.Range("B70").Value = '''45%'
$ws.Range("C70").Value = 'Comedy, Romance'
$ws.Range("A71").Value = 'Dark Waters (2019)'
$ws.Range("B71").Value = '''95%'
$ws.Range("A72").Value = 'Rosie (2019)'
$ws.Range("B72").Value = '''77%'
$ws.Range("A73").Value = 'The Mustang (2019)'
$ws.Range("B73").Value = '''74%'
$ws.Range("C73").Value = 'Drama'
$ws.Range("A74").Value = 'The Garden Left Behind (2019)'
$ws.Range("B74").ClearContents()
$ws.Range("C74").Value = 'Drama, Gay & Lesbian'
$ws.Range("A75").Value = 'Working Woman (2019)'
$ws.Range("B75").Value = '''100%'
$ws.Range("C75").Value = 'Drama'
$ws.Range("A76").Value = 'American Factory (2019)'
$ws.Range("B76").Value = '''79%'
$ws.Range("C76").Value = 'Documentary'
$ws.Range("A77").Value = 'Making Waves: The Art of Cinematic Sound (2019)'
$ws.Range("B77").Value = '''92%'
$ws.Range("A78").Value = 'Fiddler: A Miracle of Miracles (2019)'
$ws.Range("B78").Value = '''98%'
$ws.Range("C78").Value = 'Documentary'
$ws.Range("B79").Value = '''81%'
$ws.Range("A80").Value = 'The Kingmaker (2019)'
$ws.Range("B80").Value = '''90%'
$ws.Range("C80").Value = 'Documentary'
$ws.Range("A81").Value = 'I Lost My Body (2019)'
$ws.Range("B81").ClearContents()
$ws.Range("C81").Value = 'Animation, Drama, Science Fiction & Fantasy'
$ws.Range("A82").Value = 'Shadow (2019)'
$ws.Range("B82").Value = '''80%'
$ws.Range("C82").Value = 'Action & Adventure, Art House & International, Drama'
$ws.Range("A83").Value = 'Ray & Liz (2019)'
$ws.Range("B83").ClearContents()
$ws.Range("C83").Value = 'Art House & International, Drama'
$ws.Range("A84").Value = 'In Fabric (2019)'
$ws.Range("B84").Value = '''48%'
$ws.Range("C84").Value = 'Comedy, Horror'
$ws.Range("A85").Value = 'The Two Popes (2019)'
$ws.Range("B85").Value = '''88%'
$ws.Range("C85").Value = 'Drama'
$ws.Range("A86").Value = 'High Flying Bird (2019)'
$ws.Range("B86").Value = '''42%'
$ws.Range("C86").Value = 'Drama'
$ws.Range("A87").Value = 'Be Natural: The Untold Story of Alice Guy-Blaché (2019)'
$ws.Range("B87").Value = '''96%'
$ws.Range("C87").Value = 'Documentary'
$ws.Range("A88").Value = 'Love, Antosha (2019)'
$ws.Range("C88").Value = 'Documentary'
$ws.Range("A89").Value = 'Monos (2019)'
$ws.Range("B89").Value = '''85%'
$ws.Range("A90").Value = 'Harpoon (2019)'
$ws.Range("B90").Value = '''62%'
$ws.Range("C90").Value = 'Art House & International, Comedy, Horror'
$ws.Range("A92").Value = 'An Elephant Sitting Still (Da xiang xi di er zuo) (2019)'
$ws.Range("B92").Value = '''75%'
$ws.Range("C92").Value = 'Art House & International, Drama'
$ws.Range("A93").Value = 'Luce (2019)'
$ws.Range("B93").Value = '''77%'
$ws.Range("C93").Value = 'Drama, Mystery & Suspense'
$ws.Range("A94").Value = 'Mike Wallace Is Here (2019)'
$ws.Range("B94").Value = '''86%'
$ws.Range("C94").Value = 'Documentary'
$ws.Range("A95").Value = 'Jojo Rabbit (2019)'
$ws.Range("B95").Value = '''94%'
$ws.Range("C95").Value = 'Comedy, Drama'
$ws.Range("B97").Value = '''56%'
$ws.Range("A98").Value = 'Ask Dr. Ruth (2019)'
$ws.Range("B98").Value = '''95%'
$ws.Range("C98").Value = 'Documentary'
$ws.Range("A99").Value = 'The Wild Pear Tree (Ahlat agaci) (2019)'
$ws.Range("B99").Value = '''81%'
$ws.Range("C99").Value = 'Drama'
$ws.Range("A100").Value = 'Citizen K (2019)'
$ws.Range("B100").Value = '''38%'
$ws.Range("C100").Value = 'Documentary'
$ws.Range("A101").Value = 'Diane (2019)'
$ws.Range("B101").Value = '''69%'
$ws.Range("C101").Value = 'Drama'
